# Auto-generated Excel COM-interop script
# Applies crypto price/volume/hour updates and a couple of row swaps
# as described by the commit "Updated symbol list on Fri Jan  6 13:12:16 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: writes a value as literal text (matches the workbook's inline-string cells)
# without leaving a stray NumberFormat/style behind on the cell.
function Set-TextCell {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range('D2') '255.50'
Set-TextCell $ws.Range('E2') '-0.44%'
Set-TextCell $ws.Range('G2') '13'

# Row 3
Set-TextCell $ws.Range('D3') '26.50'
Set-TextCell $ws.Range('E3') '-2.45%'
Set-TextCell $ws.Range('G3') '13'

# Row 4
Set-TextCell $ws.Range('D4') '4.648'
Set-TextCell $ws.Range('E4') '-1.37%'
Set-TextCell $ws.Range('G4') '13'

# Row 5
Set-TextCell $ws.Range('D5') '0.05925'
Set-TextCell $ws.Range('E5') '0.05%'
Set-TextCell $ws.Range('G5') '13'

# Row 6
Set-TextCell $ws.Range('D6') '6.626'
Set-TextCell $ws.Range('E6') '-0.53%'
Set-TextCell $ws.Range('G6') '13'

# Row 7
Set-TextCell $ws.Range('D7') '0.8504'
Set-TextCell $ws.Range('E7') '-2.01%'
Set-TextCell $ws.Range('G7') '13'

# Row 8
Set-TextCell $ws.Range('D8') '0.9024'
Set-TextCell $ws.Range('E8') '-4.47%'
Set-TextCell $ws.Range('G8') '13'

# Row 9
Set-TextCell $ws.Range('E9') '-1.96%'
Set-TextCell $ws.Range('G9') '13'

# Row 10
Set-TextCell $ws.Range('D10') '0.04139'
Set-TextCell $ws.Range('E10') '5.77%'
Set-TextCell $ws.Range('G10') '13'

# Row 11
Set-TextCell $ws.Range('D11') '0.06998'
Set-TextCell $ws.Range('E11') '-1.41%'
Set-TextCell $ws.Range('G11') '13'

# Row 12
Set-TextCell $ws.Range('D12') '0.03032'
Set-TextCell $ws.Range('E12') '-5.66%'
Set-TextCell $ws.Range('G12') '13'

# Row 13
Set-TextCell $ws.Range('D13') '0.09085'
Set-TextCell $ws.Range('E13') '-1.84%'
Set-TextCell $ws.Range('G13') '13'

# Row 14
Set-TextCell $ws.Range('D14') '0.001532'
Set-TextCell $ws.Range('E14') '-1.01%'
Set-TextCell $ws.Range('G14') '13'

# Row 15
Set-TextCell $ws.Range('D15') '0.0006062'
Set-TextCell $ws.Range('E15') '0.55%'
Set-TextCell $ws.Range('G15') '13'

# Row 16
Set-TextCell $ws.Range('D16') '0.006049'
Set-TextCell $ws.Range('E16') '0.01%'
Set-TextCell $ws.Range('G16') '13'

# Row 17
Set-TextCell $ws.Range('D17') '3.464'
Set-TextCell $ws.Range('E17') '-1.38%'
Set-TextCell $ws.Range('G17') '13'

# Row 18
Set-TextCell $ws.Range('D18') '3.150'
Set-TextCell $ws.Range('E18') '-1.50%'
Set-TextCell $ws.Range('G18') '13'

# Row 19
Set-TextCell $ws.Range('E19') '-1.91%'
Set-TextCell $ws.Range('G19') '13'

# Row 20
Set-TextCell $ws.Range('D20') '0.3021'
Set-TextCell $ws.Range('E20') '-3.81%'
Set-TextCell $ws.Range('G20') '13'

# Row 21
Set-TextCell $ws.Range('G21') '13'

# Row 22
Set-TextCell $ws.Range('D22') '3.870'
Set-TextCell $ws.Range('E22') '-0.98%'
Set-TextCell $ws.Range('G22') '13'

# Row 23
Set-TextCell $ws.Range('D23') '0.04201'
Set-TextCell $ws.Range('E23') '-0.43%'
Set-TextCell $ws.Range('G23') '13'

# Row 24
Set-TextCell $ws.Range('D24') '0.001215'
Set-TextCell $ws.Range('E24') '-0.41%'
Set-TextCell $ws.Range('G24') '13'

# Row 25
Set-TextCell $ws.Range('D25') '0.004706'
Set-TextCell $ws.Range('E25') '9.66%'
Set-TextCell $ws.Range('G25') '13'

# Row 26
Set-TextCell $ws.Range('D26') '0.0001200'
Set-TextCell $ws.Range('E26') '0.06%'
Set-TextCell $ws.Range('G26') '13'

# Row 27
Set-TextCell $ws.Range('D27') '0.0001524'
Set-TextCell $ws.Range('E27') '1.51%'
Set-TextCell $ws.Range('G27') '13'

# Row 28
Set-TextCell $ws.Range('G28') '13'

# Row 29
Set-TextCell $ws.Range('G29') '13'

# Row 30
Set-TextCell $ws.Range('G30') '13'

# Row 31
Set-TextCell $ws.Range('G31') '13'

# Row 32
Set-TextCell $ws.Range('G32') '13'

# Row 33
Set-TextCell $ws.Range('G33') '13'

# Row 34
Set-TextCell $ws.Range('G34') '13'

# Row 35
Set-TextCell $ws.Range('G35') '13'

# Row 36
Set-TextCell $ws.Range('G36') '13'

# Row 37
Set-TextCell $ws.Range('G37') '13'

# Row 38
Set-TextCell $ws.Range('G38') '13'

# Row 39
Set-TextCell $ws.Range('G39') '13'

# Row 40
Set-TextCell $ws.Range('D40') '0.03771'
Set-TextCell $ws.Range('E40') '-1.48%'
Set-TextCell $ws.Range('G40') '13'

# Row 41
Set-TextCell $ws.Range('B41') 'BKEXToken'
Set-TextCell $ws.Range('C41') 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws.Range('D41') '0.1095'
Set-TextCell $ws.Range('E41') '-0.71%'
Set-TextCell $ws.Range('G41') '13'

# Row 42
Set-TextCell $ws.Range('B42') 'KickToken'
Set-TextCell $ws.Range('C42') 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell $ws.Range('D42') '0.003717'
Set-TextCell $ws.Range('E42') '-40.51%'
Set-TextCell $ws.Range('G42') '13'

# Row 43
Set-TextCell $ws.Range('D43') '0.002440'
Set-TextCell $ws.Range('E43') '1.72%'
Set-TextCell $ws.Range('G43') '13'

# Row 44
Set-TextCell $ws.Range('D44') '0.01446'
Set-TextCell $ws.Range('E44') '25.86%'
Set-TextCell $ws.Range('G44') '13'

# Row 45
Set-TextCell $ws.Range('D45') '0.00005153'
Set-TextCell $ws.Range('E45') '-6.29%'
Set-TextCell $ws.Range('G45') '13'

# Row 46
Set-TextCell $ws.Range('E46') '0.04%'
Set-TextCell $ws.Range('G46') '13'

# Row 47
Set-TextCell $ws.Range('D47') '0.04001'
Set-TextCell $ws.Range('G47') '13'

# Row 48
Set-TextCell $ws.Range('E48') '9,834.82%'
Set-TextCell $ws.Range('G48') '13'

# Row 49
Set-TextCell $ws.Range('D49') '0.00002101'
Set-TextCell $ws.Range('E49') '0.04%'
Set-TextCell $ws.Range('G49') '13'

# Row 50
Set-TextCell $ws.Range('D50') '0.0002001'
Set-TextCell $ws.Range('E50') '0.04%'
Set-TextCell $ws.Range('G50') '13'

# Row 51
Set-TextCell $ws.Range('G51') '13'

Write-Host "Applied $([int]121) cell updates across 50 rows"
